$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the title text in B1 (shared string "Tracker (yo2y)" -> "Tracker (yoy)")
$ws.Range("B1").Value = "Tracker (yoy)"

# New values for B2:B46 (new model output for yoy)
$newValues = @(
    3.997431822067821,
    3.507738076960831,
    -0.6219318845783861,
    1.35041188314875,
    0.5494206840561455,
    5.154642321020941,
    0.5513447351755074,
    -2.422392786878991,
    0.8459790876528261,
    0.1157408238356528,
    3.540058817858083,
    3.548719203771689,
    1.677649134169701,
    0.728128150653462,
    -1.284171144984669,
    0.6465608374731646,
    4.028908876135762,
    -0.8853013920439023,
    3.043700757466139,
    1.792946972473097,
    -1.330779030561458,
    4.293095544127379,
    0.5271650709245401,
    1.207794540244223,
    0.6875775215922753,
    2.053657220440286,
    1.385049210071831,
    2.267413823955633,
    2.181221137737266,
    0.5340766838090261,
    2.561469987863085,
    -0.8434777660186188,
    2.742529675130223,
    -0.9456210172690382,
    3.175229136339963,
    2.280434039863155,
    -3.537696405876434,
    1.864974102403472,
    -2.600824064847018,
    1.739708944864859,
    1.109473804490002,
    -0.1696800102653051,
    -0.05158516532106638,
    -0.7357271084382444,
    1.429893097472057
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
